$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column (Price) updates: force text type so Excel does not
# reinterpret numeric-looking strings (e.g. "526.35") as numbers,
# matching the original inline-string cell type. NumberFormat is
# reset back to "Normal" style afterward so no visible style diff
# is left behind on the cell.

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "57.849.90"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +1.29%  "

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "3.111.19"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +2.16%  "

$ws.Cells.Item(4, 5).Value = "  -0.02%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "526.35"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.20%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "141.07"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.30%  "

$ws.Cells.Item(7, 5).Value = "  -0.04%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "3.109.75"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +2.21%  "

$ws.Cells.Item(9, 5).Value = "  +0.16%  "

$ws.Cells.Item(10, 5).Value = "  +0.00%  "

$ws.Cells.Item(11, 5).Value = "  +2.19%  "

$ws.Cells.Item(12, 5).Value = "  +3.53%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "3.643.82"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +1.95%  "

$ws.Cells.Item(14, 5).Value = "  +1.89%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "26.28"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +4.04%  "

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "0.0000165"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +2.07%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "57.936.94"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +1.21%  "

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "3.112.81"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +2.11%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "6.12"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.24%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "12.86"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.51%  "

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "8.10"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.81%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "337.53"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +1.85%  "

$ws.Cells.Item(23, 5).Value = "  +0.02%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "0.511"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +2.85%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "66.67"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.50%  "

$ws.Cells.Item(26, 5).Value = "  +0.10%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.15%  "

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₃0932"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +4.20%  "

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "6.58"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +4.81%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "7.24"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +1.63%  "

$ws.Cells.Item(32, 5).Value = "  +3.26%  "

$ws.Cells.Item(33, 5).Value = "  +3.96%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "20.97"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +1.13%  "

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "154.03"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.24%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "4.65"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +5.77%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "6.12"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +4.21%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "26.91"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.29%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "1.31"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +3.11%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.0669"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.05%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "3.149.42"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +2.01%  "

$ws.Cells.Item(42, 5).Value = "  +5.58%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "3.91"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.41%  "

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "36.94"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.37%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "1.50"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +9.24%  "

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.00%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "2.299.21"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +2.25%  "

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "0.0260"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +2.13%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +8.35%  "

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "20.94"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +4.37%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "6.01"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +3.01%  "

